$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.793046666666667
$ws.Range("H2").Value = 11.37914
$ws.Range("I2").Value = 0.175720728331298
$ws.Range("J2").Value = 0.175720728331298
$ws.Range("M2").Value = 44.917469
$ws.Range("N2").Value = 89.83493799999999
$ws.Range("O2").Value = 0.1822888946806947
$ws.Range("P2").Value = 0.1310339777180443
$ws.Range("Q2").Value = 170.3740560655533
$ws.Range("R2").Value = 1022.24433639332
$ws.Range("S2").Value = 0.03203193733999894
$ws.Range("T2").Value = 0.02302538600076182

$ws.Range("G3").Value = 3.793046666666667
$ws.Range("H3").Value = 11.37914
$ws.Range("I3").Value = 0.175720728331298
$ws.Range("J3").Value = 0.175720728331298
$ws.Range("O3").Value = 0.02728303986213551
$ws.Range("P3").Value = 0.02941763328729693
$ws.Range("Q3").Value = 25.49975505228889
$ws.Range("R3").Value = 229.4977954706
$ws.Range("S3").Value = 0.004794195635666287
$ws.Range("T3").Value = 0.005169287947026853

$ws.Range("G4").Value = 3.793046666666667
$ws.Range("H4").Value = 11.37914
$ws.Range("I4").Value = 0.175720728331298
$ws.Range("J4").Value = 0.175720728331298
$ws.Range("M4").Value = 70.42679733333334
$ws.Range("N4").Value = 211.280392
$ws.Range("O4").Value = 0.2858135894031481
$ws.Range("P4").Value = 0.3081753134575289
$ws.Range("Q4").Value = 267.1321288692089
$ws.Range("R4").Value = 2404.18915982288
$ws.Range("S4").Value = 0.05022337209690373
$ws.Range("T4").Value = 0.05415279053448303

$ws.Range("G5").Value = 3.793046666666667
$ws.Range("H5").Value = 11.37914
$ws.Range("I5").Value = 0.175720728331298
$ws.Range("J5").Value = 0.175720728331298
$ws.Range("M5").Value = 8.721912
$ws.Range("N5").Value = 17.443824
$ws.Range("O5").Value = 0.03539619959402181
$ws.Range("P5").Value = 0.02544370482376786
$ws.Range("Q5").Value = 33.08261923856
$ws.Range("R5").Value = 198.49571543136
$ws.Range("S5").Value = 0.006219845972821506
$ws.Range("T5").Value = 0.004470986343079048

$ws.Range("G6").Value = 3.793046666666667
$ws.Range("H6").Value = 11.37914
$ws.Range("I6").Value = 0.175720728331298
$ws.Range("J6").Value = 0.175720728331298
$ws.Range("M6").Value = 86.42400633333334
$ws.Range("N6").Value = 259.272019
$ws.Range("O6").Value = 0.3507351803010249
$ws.Range("P6").Value = 0.3781762944007193
$ws.Range("Q6").Value = 327.8102891426289
$ws.Range("R6").Value = 2950.29260228366
$ws.Range("S6").Value = 0.06163144133390519
$ws.Range("T6").Value = 0.06645341388972575

$ws.Range("G7").Value = 3.793046666666667
$ws.Range("H7").Value = 11.37914
$ws.Range("I7").Value = 0.175720728331298
$ws.Range("J7").Value = 0.175720728331298
$ws.Range("M7").Value = 29.19520033333334
$ws.Range("N7").Value = 87.58560100000001
$ws.Range("O7").Value = 0.1184830961589751
$ws.Range("P7").Value = 0.1277530763126427
$ws.Range("Q7").Value = 110.7387573070156
$ws.Range("R7").Value = 996.6488157631401
$ws.Range("S7").Value = 0.02081993595200232
$ws.Range("T7").Value = 0.02244886361622146

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.537309333333334
$ws.Range("H8").Value = 13.611928
$ws.Range("I8").Value = 0.2102002350048587
$ws.Range("J8").Value = 0.2102002350048587
$ws.Range("M8").Value = 44.917469
$ws.Range("N8").Value = 89.83493799999999
$ws.Range("O8").Value = 0.1822888946806947
$ws.Range("P8").Value = 0.1310339777180443
$ws.Range("Q8").Value = 203.8044513234107
$ws.Range("R8").Value = 1222.826707940464
$ws.Range("S8").Value = 0.03831716850065797
$ws.Range("T8").Value = 0.02754337290995434

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.537309333333334
$ws.Range("H9").Value = 13.611928
$ws.Range("I9").Value = 0.2102002350048587
$ws.Range("J9").Value = 0.2102002350048587
$ws.Range("O9").Value = 0.02728303986213551
$ws.Range("P9").Value = 0.02941763328729693
$ws.Range("Q9").Value = 30.50325681812445
$ws.Range("R9").Value = 274.52931136312
$ws.Range("S9").Value = 0.005734901390667813
$ws.Range("T9").Value = 0.006183593430276571

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.537309333333334
$ws.Range("H10").Value = 13.611928
$ws.Range("I10").Value = 0.2102002350048587
$ws.Range("J10").Value = 0.2102002350048587
$ws.Range("M10").Value = 70.42679733333334
$ws.Range("N10").Value = 211.280392
$ws.Range("O10").Value = 0.2858135894031481
$ws.Range("P10").Value = 0.3081753134575289
$ws.Range("Q10").Value = 319.5481648573085
$ws.Range("R10").Value = 2875.933483715776
$ws.Range("S10").Value = 0.06007808366012393
$ws.Range("T10").Value = 0.06477852331146858

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 4.537309333333334
$ws.Range("H11").Value = 13.611928
$ws.Range("I11").Value = 0.2102002350048587
$ws.Range("J11").Value = 0.2102002350048587
$ws.Range("M11").Value = 8.721912
$ws.Range("N11").Value = 17.443824
$ws.Range("O11").Value = 0.03539619959402181
$ws.Range("P11").Value = 0.02544370482376786
$ws.Range("Q11").Value = 39.574012722112
$ws.Range("R11").Value = 237.444076332672
$ws.Range("S11").Value = 0.007440289472942271
$ws.Range("T11").Value = 0.005348272733350263

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 4.537309333333334
$ws.Range("H12").Value = 13.611928
$ws.Range("I12").Value = 0.2102002350048587
$ws.Range("J12").Value = 0.2102002350048587
$ws.Range("M12").Value = 86.42400633333334
$ws.Range("N12").Value = 259.272019
$ws.Range("O12").Value = 0.3507351803010249
$ws.Range("P12").Value = 0.3781762944007193
$ws.Range("Q12").Value = 392.1324505602925
$ws.Range("R12").Value = 3529.192055042632
$ws.Range("S12").Value = 0.07372461732374692
$ws.Range("T12").Value = 0.07949274595629785

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 4.537309333333334
$ws.Range("H13").Value = 13.611928
$ws.Range("I13").Value = 0.2102002350048587
$ws.Range("J13").Value = 0.2102002350048587
$ws.Range("M13").Value = 29.19520033333334
$ws.Range("N13").Value = 87.58560100000001
$ws.Range("O13").Value = 0.1184830961589751
$ws.Range("P13").Value = 0.1277530763126427
$ws.Range("Q13").Value = 132.4676549609698
$ws.Range("R13").Value = 1192.208894648728
$ws.Range("S13").Value = 0.02490517465671985
$ws.Range("T13").Value = 0.02685372666351114

$ws.Range("G14").Value = 13.255297
$ws.Range("H14").Value = 39.765891
$ws.Range("I14").Value = 0.6140790366638432
$ws.Range("J14").Value = 0.6140790366638433
$ws.Range("M14").Value = 44.917469
$ws.Range("N14").Value = 89.83493799999999
$ws.Range("O14").Value = 0.1822888946806947
$ws.Range("P14").Value = 0.1310339777180443
$ws.Range("Q14").Value = 595.3943920832929
$ws.Range("R14").Value = 3572.366352499757
$ws.Range("S14").Value = 0.1119397888400378
$ws.Range("T14").Value = 0.08046521880732818

$ws.Range("G15").Value = 13.255297
$ws.Range("H15").Value = 39.765891
$ws.Range("I15").Value = 0.6140790366638432
$ws.Range("J15").Value = 0.6140790366638433
$ws.Range("O15").Value = 0.02728303986213551
$ws.Range("P15").Value = 0.02941763328729693
$ws.Range("Q15").Value = 89.11222464404332
$ws.Range("R15").Value = 802.0100217963899
$ws.Range("S15").Value = 0.01675394283580141
$ws.Range("T15").Value = 0.01806475190999351

$ws.Range("G16").Value = 13.255297
$ws.Range("H16").Value = 39.765891
$ws.Range("I16").Value = 0.6140790366638432
$ws.Range("J16").Value = 0.6140790366638433
$ws.Range("M16").Value = 70.42679733333334
$ws.Range("N16").Value = 211.280392
$ws.Range("O16").Value = 0.2858135894031481
$ws.Range("P16").Value = 0.3081753134575289
$ws.Range("Q16").Value = 933.5281154121413
$ws.Range("R16").Value = 8401.753038709272
$ws.Range("S16").Value = 0.1755121336461204
$ws.Range("T16").Value = 0.1892439996115773

$ws.Range("G17").Value = 13.255297
$ws.Range("H17").Value = 39.765891
$ws.Range("I17").Value = 0.6140790366638432
$ws.Range("J17").Value = 0.6140790366638433
$ws.Range("M17").Value = 8.721912
$ws.Range("N17").Value = 17.443824
$ws.Range("O17").Value = 0.03539619959402181
$ws.Range("P17").Value = 0.02544370482376786
$ws.Range("Q17").Value = 115.611533967864
$ws.Range("R17").Value = 693.6692038071839
$ws.Range("S17").Value = 0.02173606414825803
$ws.Range("T17").Value = 0.01562444574733855

$ws.Range("G18").Value = 13.255297
$ws.Range("H18").Value = 39.765891
$ws.Range("I18").Value = 0.6140790366638432
$ws.Range("J18").Value = 0.6140790366638433
$ws.Range("M18").Value = 86.42400633333334
$ws.Range("N18").Value = 259.272019
$ws.Range("O18").Value = 0.3507351803010249
$ws.Range("P18").Value = 0.3781762944007193
$ws.Range("Q18").Value = 1145.575871878214
$ws.Range("R18").Value = 10310.18284690393
$ws.Range("S18").Value = 0.2153791216433727
$ws.Range("T18").Value = 0.2322301345546957

$ws.Range("G19").Value = 13.255297
$ws.Range("H19").Value = 39.765891
$ws.Range("I19").Value = 0.6140790366638432
$ws.Range("J19").Value = 0.6140790366638433
$ws.Range("M19").Value = 29.19520033333334
$ws.Range("N19").Value = 87.58560100000001
$ws.Range("O19").Value = 0.1184830961589751
$ws.Range("P19").Value = 0.1277530763126427
$ws.Range("Q19").Value = 386.9910513928323
$ws.Range("R19").Value = 3482.919462535491
$ws.Range("S19").Value = 0.07275798555025295
$ws.Range("T19").Value = 0.07845048603291008
